# Updated cryptos list (Price and Volume(1h) columns) as published by the
# scheduled GitHub Actions refresh job.
#
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the inline-string cells in the source file) instead of
# reinterpreting strings such as "322.00" or "572.87" as numbers, which
# would silently drop significant trailing zeros / formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.479.15"
$ws.Range("E2").Value = "'  -4.36%  "
$ws.Range("D3").Value = "'3.336.28"
$ws.Range("E3").Value = "'  -1.21%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'572.87"
$ws.Range("D6").Value = "'180.22"
$ws.Range("E6").Value = "'  -5.97%  "
$ws.Range("E7").Value = "'  +3.77%  "
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("E9").Value = "'  -3.50%  "
$ws.Range("D10").Value = "'6.62"
$ws.Range("E10").Value = "'  -2.00%  "
$ws.Range("E11").Value = "'  -3.25%  "
$ws.Range("D12").Value = "'3.916.02"
$ws.Range("E12").Value = "'  -1.23%  "
$ws.Range("D13").Value = "'0.135"
$ws.Range("E13").Value = "'  -0.51%  "
$ws.Range("D14").Value = "'26.91"
$ws.Range("E14").Value = "'  -5.79%  "
$ws.Range("D15").Value = "'66.632.56"
$ws.Range("E15").Value = "'  -4.17%  "
$ws.Range("E16").Value = "'  -2.88%  "
$ws.Range("D17").Value = "'3.326.70"
$ws.Range("E17").Value = "'  -1.34%  "
$ws.Range("D18").Value = "'438.10"
$ws.Range("E18").Value = "'  -3.06%  "
$ws.Range("E19").Value = "'  -2.41%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "'  -2.17%  "
$ws.Range("D21").Value = "'7.58"
$ws.Range("E21").Value = "'  -2.64%  "
$ws.Range("D22").Value = "'73.42"
$ws.Range("E22").Value = "'  -3.16%  "
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E24").Value = "'  -0.88%  "
$ws.Range("E25").Value = "'  -4.43%  "
$ws.Range("D26").Value = "'0.192"
$ws.Range("E26").Value = "'  +0.45%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("E27").Value = "'  -4.93%  "
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("E29").Value = "'  -3.35%  "
$ws.Range("E30").Value = "'  -2.34%  "
$ws.Range("E31").Value = "'  +0.04%  "
$ws.Range("D32").Value = "'5.26"
$ws.Range("E32").Value = "'  -5.33%  "
$ws.Range("D33").Value = "'6.76"
$ws.Range("E33").Value = "'  -3.36%  "
$ws.Range("E34").Value = "'  -5.08%  "
$ws.Range("D35").Value = "'161.32"
$ws.Range("E35").Value = "'  -2.27%  "
$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "'  -6.62%  "
$ws.Range("D37").Value = "'27.59"
$ws.Range("E37").Value = "'  +0.27%  "
$ws.Range("E38").Value = "'  -7.89%  "
$ws.Range("D39").Value = "'2.832.93"
$ws.Range("E39").Value = "'  +3.67%  "
$ws.Range("D40").Value = "'0.796"
$ws.Range("E40").Value = "'  -1.78%  "
$ws.Range("E41").Value = "'  -3.62%  "
$ws.Range("D42").Value = "'6.17"
$ws.Range("E42").Value = "'  -6.29%  "
$ws.Range("E43").Value = "'  -1.94%  "
$ws.Range("D44").Value = "'0.0666"
$ws.Range("E44").Value = "'  -3.26%  "
$ws.Range("D45").Value = "'24.28"
$ws.Range("E45").Value = "'  -4.98%  "
$ws.Range("D46").Value = "'2.32"
$ws.Range("E46").Value = "'  -8.11%  "
$ws.Range("D47").Value = "'322.00"
$ws.Range("E47").Value = "'  -4.73%  "
$ws.Range("D48").Value = "'0.0273"
$ws.Range("E48").Value = "'  -3.90%  "
$ws.Range("E49").Value = "'  +1.15%  "
$ws.Range("D50").Value = "'0.975"
$ws.Range("E50").Value = "'  -3.81%  "
$ws.Range("E51").Value = "'  -2.56%  "
